$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1 = 'winter panta for women'
    2 = 'winter recovery'
    3 = 'winter ridding pants'
    4 = 'winter runing gear men'
    5 = 'winter running'
    6 = 'winter running gear'
    7 = 'winter running gear for men'
    8 = 'winter running gear for women'
    9 = 'winter running gear men'
    10 = 'winter running gear women'
    11 = 'winter running gear womens'
    12 = 'winter running leggings for women'
    13 = 'winter running men'
    14 = 'winter running pants'
    15 = 'winter running pants women'
    16 = 'winter running tight women'
    17 = 'winter running tights for women'
    18 = 'winter running tights men'
    19 = 'winter running tights women'
    20 = 'winter running women'
    21 = 'winter skirts for women black'
    22 = 'winter sports for women'
    23 = 'winter sports leggings women'
    24 = 'winter things'
    25 = 'winter things for women'
    26 = 'winter wear for women'
    27 = 'winter weather running gear'
    28 = 'winter wishes'
    29 = 'winter women running gear'
    30 = 'winter womens running gear'
    31 = 'winter workout clothes for women'
    32 = 'winter workout pants'
    33 = 'witchy clothes for women'
    34 = 'with amazing support'
    35 = 'wolford neon 40'
    36 = 'woman addias'
    37 = 'woman apparel'
    38 = 'woman capri'
    39 = 'woman capri leggings'
    40 = 'woman capri pants'
    41 = 'woman compression pants'
    42 = 'woman compression shorts'
    43 = 'woman exercise clothes'
    44 = 'woman gym leggings'
    45 = 'woman knee high'
    46 = 'woman knee length shorts'
    47 = 'woman knee support'
    48 = 'woman legging shorts'
    49 = 'woman nike clothes'
    50 = 'woman nike tights'
    51 = 'woman pants'
    52 = 'woman rainbow'
    53 = 'woman running'
    54 = 'woman running tights'
    55 = 'woman s clothing'
    56 = 'woman skiing pants'
    57 = 'woman sport leggings'
    58 = 'woman sports wear'
    59 = 'woman tennis pants'
    60 = 'woman tight short'
    61 = 'woman tight shorts'
    62 = 'woman tights'
    63 = 'woman wear'
    64 = 'woman winter tights'
    65 = 'woman workout clothes'
    66 = 'woman workout clothing'
    67 = 'woman workout gear'
    68 = 'woman yoga'
    69 = 'womans black clothes'
    70 = 'womans black leggings capri'
    71 = 'womans black tights'
    72 = 'womans black tights control top'
    73 = 'womans capri pants'
    74 = 'womans capris'
    75 = 'womans cold weather gear'
    76 = 'womans compression'
    77 = 'womans compression capris'
    78 = 'womans compression leggings'
    79 = 'womans compression pants'
    80 = 'womans compression shorts'
    81 = 'womans compression underwear'
    82 = 'womans down pants'
    83 = 'womans gym shorts'
    84 = 'womans knee support'
    85 = 'womans legging'
    86 = 'womans legging shorts'
    87 = 'womans leggings capri'
    88 = 'womans leggings long'
    89 = 'womans leggings multi pack'
    90 = 'womans long tops for leggings'
    91 = 'womans nike pros'
    92 = 'womans pants'
    93 = 'womans pants capri'
    94 = 'womans running gear'
    95 = 'womans running leggings'
    96 = 'womans shorts clearance'
    97 = 'womans skis'
    98 = 'womans sports wear'
    99 = 'womans tight leggings'
    100 = 'womans tights'
}

foreach ($row in 1..100) {
    $ws.Range("A$row").Value = $values[$row]
}
